$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First table ("JIRA REPORT GENERATOR" banner): remove the second,
#    empty row.
# ------------------------------------------------------------------
$bannerTable = $d.Tables.Item(1)
$bannerTable.Rows.Item(2).Delete()

# ------------------------------------------------------------------
# 2) Second table (the report grid): append two new columns,
#    "SPRINT" and "EPIC LINK", with header + per-row values.
# ------------------------------------------------------------------
$reportTable = $d.Tables.Item(2)

# Add the two columns at the end of the table.
$reportTable.Columns.Add($reportTable.Columns.Item($reportTable.Columns.Count + 1))
$reportTable.Columns.Add($reportTable.Columns.Item($reportTable.Columns.Count + 1))

$sprintColIndex = $reportTable.Columns.Count - 1
$epicColIndex = $reportTable.Columns.Count

function Set-CellText($table, $row, $col, $text, $bold, $italic) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $text
    $rng = $cell.Range
    $textRange = $d.Range($rng.Start, $rng.End - 1)
    if ($bold) {
        $textRange.Font.Bold = $true
    }
    if ($italic) {
        $textRange.Font.Italic = $true
    }
}

# Header row.
Set-CellText $reportTable 1 $sprintColIndex "SPRINT" $true $true
Set-CellText $reportTable 1 $epicColIndex "EPIC LINK" $true $true

# Data rows: Row index => (Sprint, Epic Link)
$rowData = @(
    @{ Row = 2;  Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 3;  Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 4;  Sprint = "Sprint 6"; Epic = "IN-176" },
    @{ Row = 5;  Sprint = "Sprint 6"; Epic = "IN-176" },
    @{ Row = 6;  Sprint = "Sprint 6"; Epic = "IN-176" },
    @{ Row = 7;  Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 8;  Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 9;  Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 10; Sprint = "Sprint 6"; Epic = "IN-97" },
    @{ Row = 11; Sprint = "Sprint 5"; Epic = "IN-97" },
    @{ Row = 12; Sprint = "Sprint 5"; Epic = "IN-175" },
    @{ Row = 13; Sprint = "Sprint 5"; Epic = "IN-175" },
    @{ Row = 14; Sprint = "Sprint 5"; Epic = "IN-97" },
    @{ Row = 15; Sprint = "Sprint 5"; Epic = "IN-123" },
    @{ Row = 16; Sprint = "Sprint 6"; Epic = "IN-175" },
    @{ Row = 17; Sprint = "Sprint 5"; Epic = "IN-123" }
)

foreach ($entry in $rowData) {
    Set-CellText $reportTable $entry.Row $sprintColIndex $entry.Sprint $false $false
    Set-CellText $reportTable $entry.Row $epicColIndex $entry.Epic $false $false
}
